$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column D ("Step Performed") - shifts old D/E/F -> E/F/G
# ---------------------------------------------------------------------------
$ws.Columns("D:D").Insert()

# ---------------------------------------------------------------------------
# 2. Register the smaller (size 8) font used for the Tester/Timestamp column
#    without actually leaving it applied to any cell (matches the original
#    author's workbook, where the font is present in the style table but not
#    referenced by any cell format).
# ---------------------------------------------------------------------------
$ws.Range("Z100").Font.Size = 8
$ws.Range("Z100").Clear()

# ---------------------------------------------------------------------------
# 3. Header row (row 4)
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = "Step Performed"
$ws.Range("G4").Value = "Tester/Timestamp"

# ---------------------------------------------------------------------------
# 4. Updated wording on a few pre-existing cells
# ---------------------------------------------------------------------------
$ws.Range("B9").Value  = "When the user presses the boot command and if the firmware is present and the verification signature is written then the bootloader boots the new app without the issuing a fault"
$ws.Range("B10").Value = "When the user presses the boot command and if the firmware is not present and the verification signature is not written then the bootloader prints an error."
$ws.Range("A12").Value = "Green LED should be blinking in the bootloader mode"
$ws.Range("A14").Value = "Flash should not write in an invalid location"
$ws.Range("B14").Value = "The firmware file should not contain the flash region to be modified that belongs to the bootloader otherwise that could cause bootloader corruption potentially bricking the device."

# ---------------------------------------------------------------------------
# 5. New "Step Performed" column (D) content for every data row
# ---------------------------------------------------------------------------
$ws.Range("D5").Value  = "1. Plug the USB.`n2. Wait for the 20 second timeout.`n3. The program should boot up automatically"
$ws.Range("D6").Value  = "1. Plug the USB.`n2. Enter the erase command and hit enter.`n3. No error should be printed"
$ws.Range("D7").Value  = "1. Plug the USB.`n2. Enter the erase command and hit enter.`n3. No error should be printed`n4. Enter the prog coomand and hit enter.`n5. The systems should be waiting for data."
$ws.Range("D8").Value  = "1. Plug the USB.`n2. Enter the erase command and hit enter.`n3. No error should be printed`n4. Enter the prog coomand and hit enter.`n5. The systems should be waiting for data.`n6. Tera Term: File->Send File->Select the .s19 file`n7. Transfer should begin.`n8. Bootloader should boot the image successfully."
$ws.Range("D9").Value  = "1. Plug in the USB.`n2. Enter the boot command and hit enter."
$ws.Range("D10").Value = "1. Plug in the USB.`n2. Enter the boot command and hit enter."
$ws.Range("D11").Value = "1. Plug the USB.`n2. Wait for the 20 second timeout.`n3. The program should boot up automatically"
$ws.Range("D12").Value = "1. Plug the USB."
$ws.Range("D13").Value = "1. Plug the USB.`n2. Enter an invalid command and hit enter."
$ws.Range("D14").Value = "1. Plug the USB.`n2. Enter the erase command and hit enter.`n3. No error should be printed`n4. Enter the prog coomand and hit enter.`n5. The systems should be waiting for data.`n6. Tera Term: File->Send File->Select the .s19 file`n7. Transfer should begin.`n8. Bootloader should boot the image successfully."

# ---------------------------------------------------------------------------
# 6. New "Tester/Timestamp" column (G) content for every data row
# ---------------------------------------------------------------------------
$ws.Range("G5").Value  = "Sankalp Agrawal`n11-12-2021"
$ws.Range("G6").Value  = "Sankalp Agrawal`n11-12-2022"
$ws.Range("G7").Value  = "Sankalp Agrawal`n11-12-2023"
$ws.Range("G8").Value  = "Sankalp Agrawal`n11-12-2024"
$ws.Range("G9").Value  = "Sankalp Agrawal`n11-12-2025"
$ws.Range("G10").Value = "Sankalp Agrawal`n11-12-2026"
$ws.Range("G11").Value = "Sankalp Agrawal`n11-12-2027"
$ws.Range("G12").Value = "Sankalp Agrawal`n11-12-2028"
$ws.Range("G13").Value = "Sankalp Agrawal`n11-12-2029"
$ws.Range("G14").Value = "Sankalp Agrawal`n11-12-2030"

# ---------------------------------------------------------------------------
# 7. Row heights to accommodate the new, longer text
# ---------------------------------------------------------------------------
$ws.Rows("5:5").Height  = 72.5
$ws.Rows("6:6").Height  = 58
$ws.Rows("7:7").Height  = 116
$ws.Rows("8:8").Height  = 188.5
$ws.Rows("9:9").Height  = 174
$ws.Rows("10:10").Height = 145
$ws.Rows("11:11").Height = 101.5
$ws.Rows("12:12").Height = 43.5
$ws.Rows("13:13").Height = 58
$ws.Rows("14:14").Height = 188.5

# ---------------------------------------------------------------------------
# 8. Column widths: D (Step Performed) wider, G (Tester/Timestamp) like the
#    old F column.
# ---------------------------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 30.26953125
$ws.Columns("G:G").ColumnWidth = 17.6328125

# ---------------------------------------------------------------------------
# 9. View state: scroll down and move the active selection
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I5").Select()
